$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89").Value = "2024-10-30 00:00:00"
$ws.Range("B89").Value = 73850
$ws.Range("C89").Value = 10338.07
$ws.Range("D89").Value = 9148.73
$ws.Range("E89").Value = 7.1216
